# Assignment 4: Ready for hand-in
# Fill in the remaining "A-Grade" test-summary rows and mark the previously
# empty/"Bad"-styled D2 cell with its real answer.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("A-Grade")

# --- New rows describing further unit tests (entered first, a row at a
# time for the Method Name / Purpose / Test result columns, then the Action
# column for all three rows) --------------------------------------------
$ws.Range("A4").Value = "BugViewModel.ShowClosedReason"
$ws.Range("B4").Value = "If status changed to Finished or Rejected then ShowCloseReason should be set to true"
$ws.Range("C4").Value = "It was as expected"

$ws.Range("A5").Value = "BugViewModel.Bug.CloseReason"
$ws.Range("B5").Value = "If status changed from Finished or Rejected then Bug.CloseReason is set to empty string"
$ws.Range("C5").Value = "It was as expected"

$ws.Range("A6").Value = "MainViewModel.Bugs.CollectionChanged"
$ws.Range("B6").Value = "When colletion change on Add/Delete the text indiciating how many bugs there are in the system should change"
$ws.Range("C6").Value = "It was as expected"

$ws.Range("D4").Value = "None"
$ws.Range("D5").Value = "None"
$ws.Range("D6").Value = "None"

# Match formatting (wrap text) used by the other rows in the table, and set
# the row heights to match the content.
$ws.Range("B2").Copy()
$ws.Range("B4:B6").PasteSpecial(-4122)   # xlPasteFormats

$ws.Rows.Item(4).RowHeight = 60
$ws.Rows.Item(5).RowHeight = 60
$ws.Rows.Item(6).RowHeight = 75

# --- D2: was an empty cell styled with the built-in "Bad" cell style ------
# Give it the same (wrap-text) formatting used by the rest of the table,
# then fill in the answer text (entered last).
$ws.Range("C2").Copy()
$ws.Range("D2").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("D2").Value = "I made validate and save public methods"

# The "Bad" cell style is no longer used anywhere in the workbook - remove it.
$wb.Styles.Item("Dålig").Delete()

# Move the selection to D2, mirroring the author's final cursor position.
$ws.Range("D2").Select()
